# Updates the "cryptos" price/volume snapshot to the values captured in the
# latest GitHub Actions run. Numeric-looking price strings in column D are
# forced back to text (they were stored as text in the source data, e.g.
# "72.378.16" / "1.00") so Excel's automatic number coercion doesn't change
# their stored type; the temporary "@" number format is reset to "Normal"
# right after the write so no residual style index is left on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '72.378.16'
$ws.Range("E2").Value = '  +4.17%  '
$ws.Range("D3").Value = '2.623.09'
$ws.Range("E3").Value = '  +4.35%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '604.19'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.19%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '178.53'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.29%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.525'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.21%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.172'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +9.44%  '
$ws.Range("D10").Value = '2.621.56'
$ws.Range("E10").Value = '  +4.31%  '
$ws.Range("E11").Value = '  +1.06%  '
$ws.Range("E12").Value = '  +2.60%  '
$ws.Range("E13").Value = '  +0.52%  '
$ws.Range("D14").Value = '3.106.81'
$ws.Range("E14").Value = '  +4.34%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000187'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.68%  '
$ws.Range("D16").Value = '72.162.32'
$ws.Range("E16").Value = '  +4.03%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.54'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.35%  '
$ws.Range("D18").Value = '2.621.83'
$ws.Range("E18").Value = '  +3.97%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '381.59'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.40%  '
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.93'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.81%  '
$ws.Range("B21").Value = 'Chainlink'
$ws.Range("C21").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.56'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.88%  '
$ws.Range("E22").Value = '  +1.89%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.02'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +17.62%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.09'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.68%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.37'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.22%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.91'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +8.89%  '
$ws.Range("D28").Value = '2.752.69'
$ws.Range("E28").Value = '  +5.62%  '
$ws.Range("E29").Value = '  +0.08%  '
$ws.Range("D30").Value = '0.0₃0950'
$ws.Range("E30").Value = '  +6.11%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '519.14'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.03%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.03'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.02%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.33'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.85%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.83'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.63%  '
$ws.Range("E35").Value = '  -0.11%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '163.68'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.30%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '19.26'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.90%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.09'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.02%  '
$ws.Range("E39").Value = '  +6.11%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.112'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.37%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.83'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.50%  '
$ws.Range("E42").Value = '  -0.01%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.04'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.21%  '
$ws.Range("E44").Value = '  +9.07%  '
$ws.Range("E45").Value = '  +3.21%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '39.49'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.79%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '150.18'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.21%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.68'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.76%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.542'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.06%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.70'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +7.46%  '
$ws.Range("E51").Value = '  +3.94%  '

$wb.Save()
